# Fruta / hortaliza, semanal
# Weekly data refresh: insert a new latest-week record for
# "Vega Central Mapocho de Santiago - Frambuesa" at row 67, pushing the
# previously-existing rows (67-108) down by one (they become 68-109) while
# keeping all of their data intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 67; Excel shifts rows 67:108 down to 68:109,
# carrying their values/styles with them and growing the used range to T109.
$ws.Range("A67:T67").EntireRow.Insert()

# Populate the newly inserted row with this week's record.
$ws.Range("A67").Value = 9
$ws.Range("B67").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C67").Value = "Metropolitana"
$ws.Range("D67").Value = 44907
$ws.Range("E67").Value = 13
$ws.Range("F67").Value = "Fruta"
$ws.Range("G67").Value = 100101
$ws.Range("H67").Value = "Berries"
$ws.Range("I67").Value = 100101004
$ws.Range("J67").Value = "Frambuesa"
$ws.Range("K67").Value = "Sin especificar"
$ws.Range("L67").Value = "Primera"
$ws.Range("M67").Value = 470
$ws.Range("N67").Value = 8000
$ws.Range("O67").Value = 8500
$ws.Range("P67").Value = 8266
$ws.Range("Q67").Value = "`$/bandeja 2 kilos"
$ws.Range("R67").Value = "Provincia de Curicó"
$ws.Range("S67").Value = 4133
$ws.Range("T67").Value = 2
